$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected against edits to locked cells (accidental-edit
# guard, no password needed to lift it). Temporarily unprotect, make the
# data updates, then restore protection.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer banner (A59).
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns for each holding
# row (2-55) plus the Total row's Percent Change (E56).
$ws.Range("D2").Value = 0.02071974535328139
$ws.Range("E2").Value = -0.0215851602023609
$ws.Range("D3").Value = 0.01852647719904582
$ws.Range("E3").Value = -0.00387972841901052
$ws.Range("D4").Value = 0.01938781163849709
$ws.Range("E4").Value = -0.002471576866040492
$ws.Range("D5").Value = 0.02064187784484699
$ws.Range("E5").Value = -0.013793103448276
$ws.Range("D6").Value = 0.01967532243887026
$ws.Range("E6").Value = -0.02338041889917186
$ws.Range("D7").Value = 0.01949842343252954
$ws.Range("E7").Value = 0.007905138339920903
$ws.Range("D8").Value = 0.01997401421481349
$ws.Range("E8").Value = -0.002459016393442703
$ws.Range("D9").Value = 0.02036055651309297
$ws.Range("E9").Value = 0.002039695614654446
$ws.Range("D10").Value = 0.01901883944468486
$ws.Range("E10").Value = -0.005165028974552732
$ws.Range("D11").Value = 0.02014971525948598
$ws.Range("E11").Value = -0.008174791914387636
$ws.Range("D12").Value = 0.01918835102073821
$ws.Range("E12").Value = -0.006711409395973256
$ws.Range("D13").Value = 0.02096952036110557
$ws.Range("E13").Value = 0.001942376173518889
$ws.Range("D14").Value = 0.01987697931968755
$ws.Range("E14").Value = -0.001968780762199418
$ws.Range("D15").Value = 0.01889664735452627
$ws.Range("E15").Value = -0.0004649000464900244
$ws.Range("D16").Value = 0.01782646826424838
$ws.Range("E16").Value = 0.004659289458357474
$ws.Range("D17").Value = 0.0185144975823636
$ws.Range("E17").Value = 0.02431791221826818
$ws.Range("D18").Value = 0.01605428363639267
$ws.Range("E18").Value = 0.02631578947368407
$ws.Range("D19").Value = 0.01487528969458466
$ws.Range("E19").Value = 0.005019932083271783
$ws.Range("D20").Value = 0.02212635201205149
$ws.Range("E20").Value = 0.002346146904890878
$ws.Range("D21").Value = 0.02171065931317862
$ws.Range("E21").Value = 0.01771229928819729
$ws.Range("D22").Value = 0.02112825028214493
$ws.Range("E22").Value = 0.01682085786375098
$ws.Range("D23").Value = 0.02064387444762736
$ws.Range("E23").Value = 0.001015522994342044
$ws.Range("D24").Value = 0.01875808312156865
$ws.Range("E24").Value = 0.01852048962213937
$ws.Range("D25").Value = 0.01904279867804929
$ws.Range("E25").Value = 0.01895655091427617
$ws.Range("D26").Value = 0.02002392928432272
$ws.Range("E26").Value = -0.0003589590188455283
$ws.Range("D27").Value = 0.01880959547330217
$ws.Range("E27").Value = -0.00830078125
$ws.Range("D28").Value = 0.02010239577359123
$ws.Range("E28").Value = 0.02523762700753851
$ws.Range("D29").Value = 0.01693917798865231
$ws.Range("E29").Value = 0.0180339462517678
$ws.Range("D30").Value = 0.01307655024975005
$ws.Range("E30").Value = -0.02759031361651454
$ws.Range("D31").Value = 0.009811306062734254
$ws.Range("E31").Value = 0.007122507122507171
$ws.Range("D32").Value = 0.01670357886056874
$ws.Range("E32").Value = 0.0372340425531914
$ws.Range("D33").Value = 0.0190785378678179
$ws.Range("E33").Value = -0.004248861911987922
$ws.Range("D34").Value = 0.02022398888291571
$ws.Range("E34").Value = -0.06948228882833785
$ws.Range("D35").Value = 0.01728658687243654
$ws.Range("E35").Value = 0.008547008547008739
$ws.Range("D36").Value = 0.02002912045155168
$ws.Range("E36").Value = 0.0111248454882571
$ws.Range("D37").Value = 0.01763639167955723
$ws.Range("E37").Value = 0
$ws.Range("D38").Value = 0.02072174195606176
$ws.Range("E38").Value = -0.02384737678855331
$ws.Range("D39").Value = 0.02295933469202152
$ws.Range("E39").Value = 0.004313343536941794
$ws.Range("D40").Value = 0.01949882275308562
$ws.Range("E40").Value = 0.01478599221789878
$ws.Range("D41").Value = 0.02107773623180158
$ws.Range("E41").Value = -0.02243103970900284
$ws.Range("D42").Value = 0.0196998806530688
$ws.Range("E42").Value = -0.01193914885422676
$ws.Range("D43").Value = 0.01998799043427607
$ws.Range("E43").Value = -0.008510638297872353
$ws.Range("D44").Value = 0.01959026716042653
$ws.Range("E44").Value = -0.0004688232536332881
$ws.Range("D45").Value = 0.01917736970544618
$ws.Range("E45").Value = -0.01795939614783981
$ws.Range("D46").Value = 0.01953356364146404
$ws.Range("E46").Value = -0.009219698673262777
$ws.Range("D47").Value = 0.01954294767453178
$ws.Range("E47").Value = -0.01458914396052347
$ws.Range("D48").Value = 0.01858038547411579
$ws.Range("E48").Value = 0.02385557704706653
$ws.Range("D49").Value = 0.01749623016437531
$ws.Range("E49").Value = 0.02047244094488199
$ws.Range("D50").Value = 0.01791911063325751
$ws.Range("E50").Value = -0.01163257119935823
$ws.Range("D51").Value = 0.01734219225986983
$ws.Range("E51").Value = 0.02430965305640775
$ws.Range("D52").Value = 0.0178017103897718
$ws.Range("E52").Value = 0.02691790040376874
$ws.Range("D53").Value = 0.01688287379024589
$ws.Range("E53").Value = 0.002814636107760338
$ws.Range("D54").Value = 0.007463301193020075
$ws.Range("E54").Value = -0.01966292134831471
$ws.Range("D55").Value = 0.007438543318543497
$ws.Range("E55").Value = -0.008803951041443048
$ws.Range("E56").Value = 0.0003776574159064516

# Restore sheet protection.
$ws.Protect()
